$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.924.11'
$ws.Range("E2").Value = '  -1.75%  '

$ws.Range("D3").Value = '1.828.26'
$ws.Range("E3").Value = '  -2.25%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.04'
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6861'
$ws.Range("E6").Value = '  -2.95%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07617'
$ws.Range("E8").Value = '  -3.33%  '

$ws.Range("E9").Value = '  -4.46%  '

$ws.Range("E10").Value = '  -4.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07753'
$ws.Range("E11").Value = '  -3.19%  '

$ws.Range("D12").Value = '1.834.83'
$ws.Range("E12").Value = '  -2.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.046'
$ws.Range("E13").Value = '  -3.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.20'
$ws.Range("E14").Value = '  -4.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6732'
$ws.Range("E15").Value = '  -4.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.433'

$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").Value = '28.942.43'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.12'
$ws.Range("E19").Value = '  -5.74%  '

$ws.Range("D20").Value = '2.098.78'
$ws.Range("E20").Value = '  -1.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.65'
$ws.Range("E21").Value = '  -4.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.414'
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1472'
$ws.Range("E25").Value = '  -5.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.27'
$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.719'
$ws.Range("E27").Value = '  -3.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.14'
$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.537'
$ws.Range("E29").Value = '  +2.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.202'
$ws.Range("E30").Value = '  -3.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.157'
$ws.Range("E31").Value = '  -2.40%  '

$ws.Range("E32").Value = '  -1.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05111'
$ws.Range("E33").Value = '  -4.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7638'
$ws.Range("E34").Value = '  +2.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.816'
$ws.Range("E35").Value = '  -4.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.143'
$ws.Range("E36").Value = '  -2.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.697'
$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01834'
$ws.Range("E38").Value = '  -2.58%  '

$ws.Range("D39").Value = '1.220.69'
$ws.Range("E39").Value = '  -3.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.703'
$ws.Range("E40").Value = '  -1.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9150'
$ws.Range("E41").Value = '  +1.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.05'
$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9993'
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").Value = '1.997.63'
$ws.Range("E44").Value = '  -1.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5168'
$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("E46").Value = '  -6.27%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.367'
$ws.Range("E47").Value = '  -10.16%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.491'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.01'
$ws.Range("E49").Value = '  -12.43%  '

$ws.Range("E50").Value = '  -4.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.890'
$ws.Range("E51").Value = '  -2.63%  '
